# Update the cryptos symbol-list worksheet with the latest scraped values
# (GitHub Actions refresh run). All Price (column D) and Volume(1h) (column E)
# cells are stored as text in this workbook, so we force a Text number format
# before writing any numeric-looking value to keep it from being re-interpreted
# as a number (which would drop significant digits / trailing zeros).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($address, $value) {
    $rng = $ws.Range($address)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# Row 2 - BNB
Set-TextValue "D2" "247.92"

# Row 3 - OKB
Set-TextValue "D3" "22.72"

# Row 4 - HuobiToken
Set-TextValue "D4" "5.296"

# Row 5 - Cronos
Set-TextValue "D5" "0.05733"

# Row 6 - GateToken
Set-TextValue "D6" "3.436"

# Row 7 - MXToken
Set-TextValue "D7" "0.8096"

# Row 8 - FTXToken
Set-TextValue "D8" "0.8756"

# Row 12 - BitrueCoin
Set-TextValue "D12" "0.03116"

# Row 13 - BitMartToken
Set-TextValue "D13" "0.09399"

# Row 14 - MCDex
Set-TextValue "D14" "3.890"

# Row 15 - BitForexToken
Set-TextValue "D15" "0.001577"

# Row 16 - CoinExToken
Set-TextValue "D16" "0.04817"

# Row 17 - One
Set-TextValue "D17" "0.0005852"
$ws.Range("E17").Value = "16OneONEWorstin24h"

# Row 18 - TigerCash
Set-TextValue "D18" "0.006144"

# Row 20 - BitKan
Set-TextValue "D20" "0.0009968"

# Row 23 - KuCoinToken
Set-TextValue "D23" "6.341"

# Row 24 - BTSEToken
Set-TextValue "D24" "2.190"

# Row 25 - BitpandaEcosystemToken
Set-TextValue "D25" "0.3277"

# Row 40 - IDEX
Set-TextValue "D40" "0.03947"

# Row 41 - KickToken
Set-TextValue "D41" "0.006737"

# Row 44 - LocalTraders
Set-TextValue "D44" "0.007273"

# Row 45 - CoinLion
Set-TextValue "D45" "0.00005612"

# Row 48 - BOLO
Set-TextValue "D48" "0.1776"
$ws.Range("E48").Value = "47BOLOBOLO"

# Row 49 - CryptobidCoin
Set-TextValue "D49" "0.00002101"

Write-Host "Updated symbol list with refreshed coinranking.com values"
